$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 457, shifting existing rows (457-509) down to (458-510)
$ws.Rows.Item(457).Insert()

# Populate the newly inserted row 457 with its data.
# Most fields carry over from the row that used to occupy position 457
# (now shifted to 458), except D, J, K, L, M, O and P which hold new values.
$ws.Cells.Item(457, 1).Value = 3
$ws.Cells.Item(457, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(457, 3).Value = "Coquimbo"
$ws.Cells.Item(457, 4).Value = 44946
$ws.Cells.Item(457, 5).Value = 5
$ws.Cells.Item(457, 6).Value = 100112031
$ws.Cells.Item(457, 7).Value = "Poroto verde"
$ws.Cells.Item(457, 8).Value = "Magnum"
$ws.Cells.Item(457, 9).Value = "Primera"
$ws.Cells.Item(457, 10).Value = 68
$ws.Cells.Item(457, 11).Value = 26000
$ws.Cells.Item(457, 12).Value = 27000
$ws.Cells.Item(457, 13).Value = 26559
$ws.Cells.Item(457, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(457, 15).Value = "Provincia de Santiago"
$ws.Cells.Item(457, 16).Value = 1062
$ws.Cells.Item(457, 17).Value = 25
$ws.Cells.Item(457, 18).Value = "Hortaliza"
